$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.649.73"
$ws.Range("D3").Value = "1.642.05"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'215.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("D6").Value = "'0.505"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.18%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").Value = "'0.0627"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("D10").Value = "'19.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.30%  "
$ws.Range("D11").Value = "'0.0842"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "1.871.17"
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.43%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.627.76"
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("D15").Value = "'0.530"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").Value = "'65.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("D17").Value = "26.698.53"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "0.0₃0747"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").Value = "'216.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").Value = "'1.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").Value = "'6.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.38%  "
$ws.Range("D23").Value = "'9.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("D24").Value = "'2.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +12.98%  "
$ws.Range("D25").Value = "'145.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("D28").Value = "'7.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.43%  "
$ws.Range("D29").Value = "'15.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.62%  "
$ws.Range("D30").Value = "'0.0516"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.31%  "
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("E32").Value = "  +2.64%  "
$ws.Range("E33").Value = "  +2.30%  "
$ws.Range("D34").Value = "1.280.34"
$ws.Range("E34").Value = "  +5.11%  "
$ws.Range("D35").Value = "'1.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.61%  "
$ws.Range("D36").Value = "'0.0182"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.17%  "
$ws.Range("D37").Value = "'2.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.45%  "
$ws.Range("D38").Value = "'0.532"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.30%  "
$ws.Range("D39").Value = "'0.828"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.89%  "
$ws.Range("E40").Value = "  +0.40%  "
$ws.Range("D41").Value = "'0.816"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.87%  "
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("D43").Value = "'5.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.54%  "
$ws.Range("D44").Value = "1.781.55"
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("D45").Value = "'92.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("D46").Value = "'59.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.52%  "
$ws.Range("E47").Value = "  +2.54%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.0516"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'7.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.58%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.0971"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.99%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.407"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.37%  "
